$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Row 4: TEXT ID moves from SingleUseId1 -> SingleUseId2, ALIGNMENT Right -> Left
$ws.Range("B4").Value = "SingleUseId2"
$ws.Range("D4").Value = "Left"

# Row 5: TEXT ID moves from SingleUseId2 -> SingleUseId3 (new), ALIGNMENT Left -> Right
$ws.Range("B5").Value = "SingleUseId3"
$ws.Range("D5").Value = "Right"
